$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.314.91'
$ws.Range('E2').Value = '  +1.25%  '
$ws.Range('D3').Value = '1.854.12'
$ws.Range('E3').Value = '  +1.46%  '
$ws.Range('E4').Value = '  -0.83%  '
$ws.Range('D5').Value = '''314.16'
$ws.Range('E6').Value = '  -0.68%  '
$ws.Range('D7').Value = '''0.4610'
$ws.Range('E7').Value = '  -1.01%  '
$ws.Range('D8').Value = '''0.3711'
$ws.Range('E8').Value = '  +0.20%  '
$ws.Range('D9').Value = '''0.07305'
$ws.Range('E9').Value = '  -0.82%  '
$ws.Range('D10').Value = '''0.8860'
$ws.Range('E10').Value = '  +1.48%  '
$ws.Range('E11').Value = '  +0.06%  '
$ws.Range('D12').Value = '''0.07786'
$ws.Range('E12').Value = '  -0.50%  '
$ws.Range('D13').Value = '1.823.84'
$ws.Range('E13').Value = '  -0.70%  '
$ws.Range('D14').Value = '''5.375'
$ws.Range('E14').Value = '  +0.62%  '
$ws.Range('D15').Value = '''6.551'
$ws.Range('E15').Value = '  -0.71%  '
$ws.Range('D16').Value = '''91.70'
$ws.Range('E16').Value = '  -0.02%  '
$ws.Range('D17').Value = '''1.001'
$ws.Range('E17').Value = '  -0.85%  '
$ws.Range('D18').Value = '''0.000008980'
$ws.Range('E18').Value = '  +1.63%  '
$ws.Range('E19').Value = '  -0.75%  '
$ws.Range('E20').Value = '  +1.36%  '
$ws.Range('D21').Value = '27.322.52'
$ws.Range('E21').Value = '  +0.78%  '
$ws.Range('D22').Value = '''5.130'
$ws.Range('E22').Value = '  -0.02%  '
$ws.Range('E23').Value = '  -0.52%  '
$ws.Range('D24').Value = '2.056.01'
$ws.Range('E24').Value = '  +5.23%  '
$ws.Range('D25').Value = '''1.926'
$ws.Range('E25').Value = '  +5.30%  '
$ws.Range('D26').Value = '''151.54'
$ws.Range('E26').Value = '  -0.49%  '
$ws.Range('D27').Value = '''18.40'
$ws.Range('E27').Value = '  +0.60%  '
$ws.Range('D28').Value = '''2.061'
$ws.Range('E28').Value = '  -0.82%  '
$ws.Range('B29').Value = 'BitcoinCash'
$ws.Range('C29').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D29').Value = '''116.22'
$ws.Range('E29').Value = '  +0.85%  '
$ws.Range('B30').Value = 'InternetComputer(DFINITY)'
$ws.Range('C30').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D30').Value = '''5.095'
$ws.Range('E30').Value = '  -0.24%  '
$ws.Range('D31').Value = '''0.08834'
$ws.Range('E31').Value = '  -0.51%  '
$ws.Range('D32').Value = '''3.129'
$ws.Range('E32').Value = '  +5.61%  '
$ws.Range('D33').Value = '''0.7766'
$ws.Range('E33').Value = '  +7.12%  '
$ws.Range('D34').Value = '''1.170'
$ws.Range('E34').Value = '  +3.31%  '
$ws.Range('D35').Value = '''4.499'
$ws.Range('E35').Value = '  +1.44%  '
$ws.Range('D36').Value = '''2.664'
$ws.Range('E36').Value = '  +7.36%  '
$ws.Range('D37').Value = '''0.01960'
$ws.Range('E37').Value = '  +0.43%  '
$ws.Range('D38').Value = '''1.077'
$ws.Range('E38').Value = '  +0.50%  '
$ws.Range('D39').Value = '''0.05235'
$ws.Range('E39').Value = '  +0.40%  '
$ws.Range('D40').Value = '''2.967'
$ws.Range('E40').Value = '  +1.61%  '
$ws.Range('D41').Value = '''6.992'
$ws.Range('E41').Value = '  -2.29%  '
$ws.Range('E42').Value = '  -0.76%  '
$ws.Range('E43').Value = '  +0.34%  '
$ws.Range('D44').Value = '''8.417'
$ws.Range('E44').Value = '  +2.62%  '
$ws.Range('D45').Value = '''0.4813'
$ws.Range('E45').Value = '  -0.01%  '
$ws.Range('D46').Value = '''10.24'
$ws.Range('E46').Value = '  +0.82%  '
$ws.Range('E47').Value = '  -0.79%  '
$ws.Range('D48').Value = '''1.650'
$ws.Range('E48').Value = '  +1.59%  '
$ws.Range('D49').Value = '''102.57'
$ws.Range('E49').Value = '  +0.30%  '
$ws.Range('D50').Value = '''0.06219'
$ws.Range('E50').Value = '  +0.10%  '
$ws.Range('D51').Value = '''65.56'
$ws.Range('E51').Value = '  +2.04%  '
